# Daily attendance processing - 2025-09-29 22:04:00
# Re-sync of recorded sessions: the LMS re-check found that the sessions
# previously marked "Recorded" had not actually been recorded, so they flip
# back to "Not Recorded" (losing their "Recorded By" emails), and the
# Class/Group statistics + the "Recorded By" column width are refreshed to
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Sessions that flip from "Recorded" (green) to "Not Recorded" (pink)
# ---------------------------------------------------------------------
$flippedRows = @(7, 12, 15, 24, 29, 32, 72, 112, 115, 129, 132, 141)

# A "Not Recorded" formatted row already on the sheet (row 2) - copy its
# look (fill/font/alignment) onto each flipped row instead of re-deriving
# the format from scratch, so the workbook keeps reusing the same style.
$formatSource = $ws.Range("A2:I2")

foreach ($r in $flippedRows) {
    $formatSource.Copy()
    $target = $ws.Range("A" + $r + ":I" + $r)
    $target.PasteSpecial(-4122)  # xlPasteFormats

    # "Recorded By" is cleared - nobody recorded the session after all.
    $ws.Range("G" + $r).ClearContents()

    # Status goes back to "Not Recorded".
    $ws.Range("I" + $r).Value = "Not Recorded"
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) "Recorded By" column (G) shrinks now that it no longer holds long
#    lists of e-mail addresses.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 12.17

# ---------------------------------------------------------------------
# 3) Top "Class Statistics" block
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 0     # Recorded Sessions
$ws.Range("L7").Value = 58    # Missing Sessions

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "0.0%"   # Coverage %

# ---------------------------------------------------------------------
# 4) "Group Statistics" block - Recorded/Missing/Coverage % per group,
#    for every group that had a flipped session.
# ---------------------------------------------------------------------
$groupStatRows = @(15, 16, 19, 21, 22, 23)
$recordedVal = @{15 = 0; 16 = 0; 19 = 0; 21 = 0; 22 = 0; 23 = 0}
$missingVal  = @{15 = 10; 16 = 9; 19 = 7; 21 = 7; 22 = 7; 23 = 1}

foreach ($r in $groupStatRows) {
    $ws.Range("O" + $r).Value = $recordedVal[$r]
    $ws.Range("P" + $r).Value = $missingVal[$r]

    $rCell = $ws.Range("R" + $r)
    $rCell.NumberFormat = "@"
    $rCell.Value = "0.0%"
}

Write-Output "Daily attendance processing complete."
